$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was collected; insert it as row 12 and push the
# existing historical rows (old 12-25) down to rows 13-26.
$ws.Rows("12:12").Insert()

$ws.Range("A12").Value = 1
$ws.Range("B12").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C12").Value = 'Arica y Parinacota'
$ws.Range("D12").Value = 45272
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = 100112017
$ws.Range("G12").Value = 'Ramas de apio'
$ws.Range("H12").Value = 'Sin especificar'
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 230
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 6000
$ws.Range("M12").Value = 5348
$ws.Range("N12").Value = '$/atado 7 kilos'
$ws.Range("O12").Value = 'Región de Arica y Parinacota'
$ws.Range("P12").Value = 5348
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = 'Hortaliza'
